$d = $word.ActiveDocument

# The document currently ends with an empty paragraph right before the
# section properties (sectPr). We need to insert new content (a heading,
# a description paragraph, and a table) right AFTER that empty paragraph.
#
# Simply inserting at the very end of the document content merges new
# content into that existing trailing empty paragraph, so instead we
# first create a temporary paragraph after it (which correctly lands
# after the existing paragraph mark), insert our new OOXML content at
# the start of that temporary paragraph, and finally remove the now
# redundant temporary empty paragraph that is left dangling at the end.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$anchorPara = $d.Paragraphs.Last

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Nagwek1"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:ind w:left="426" w:hanging="426"/>
            </w:pPr>
            <w:r>
              <w:t>Implementacja co najmniej jednej listy kontroli dostępu ACL wewnątrz zabezpieczonej sieci</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Lista ACL została skonfigurowana na Routerze R1. Póki co jest to tylko lista ACL blokująca ruch do VLAN30, jednak nie jestem pewien czy nie trzeba tu stworzyć ZBF opartego na strefach </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>internal</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> i </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>external</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> tak jak w instrukcji na labach.</w:t>
            </w:r>
          </w:p>
          <w:tbl>
            <w:tblPr>
              <w:tblStyle w:val="Tabela-Siatka"/>
              <w:tblW w:w="0" w:type="auto"/>
              <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
            </w:tblPr>
            <w:tblGrid>
              <w:gridCol w:w="9062"/>
            </w:tblGrid>
            <w:tr>
              <w:tc>
                <w:tcPr>
                  <w:tcW w:w="9062" w:type="dxa"/>
                </w:tcPr>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>enable</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>configure terminal</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t xml:space="preserve">access-list 100 deny </w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>ip</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t xml:space="preserve"> any 192.168.1.0 0.0.0.255</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t xml:space="preserve">access-list 100 permit </w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>ip</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t xml:space="preserve"> any </w:t>
                  </w:r>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>any</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>interface fa0/0</w:t>
                  </w:r>
                </w:p>
                <w:p>
                  <w:pPr>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                  </w:pPr>
                  <w:proofErr w:type="spellStart"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t>ip</w:t>
                  </w:r>
                  <w:proofErr w:type="spellEnd"/>
                  <w:r>
                    <w:rPr>
                      <w:lang w:val="en-US"/>
                    </w:rPr>
                    <w:t xml:space="preserve"> access-group 100 in</w:t>
                  </w:r>
                </w:p>
              </w:tc>
            </w:tr>
          </w:tbl>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)
$insertionPoint.InsertXML($xml)

# Remove the now-redundant empty helper paragraph left dangling at the
# very end of the document (right before the section properties).
$trailingPara = $d.Paragraphs.Last
$trailingPara.Range.Delete()

Write-Host "Edit applied."
